$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename ingredient "Thickener" (row 10) -> "Stabilizer".
#    This text already exists elsewhere (row 12), so the shared-string table
#    will naturally de-duplicate it, exactly like the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Stabilizer"

# ---------------------------------------------------------------------------
# 2) Bugfix: remove the 0.9 / 0.95 multipliers from the PAC/POD formulas
#    on rows 15 and 16.
# ---------------------------------------------------------------------------
$ws.Range("M15").Formula = "=((24.5*100)+(0*100)+(0.8*190))/25.3"
$ws.Range("N15").Formula = "=((24.5*100)+(0*16)+(0.8*130))/25.3"
$ws.Range("M16").Formula = "=((40*100)+(8*100) + (1.6*190))/49.6"
$ws.Range("N16").Formula = "=((40*100)+(8*16)+(1.6*130))/49.6"

# ---------------------------------------------------------------------------
# 3) Row 17 M/N values changed (bugfix to the Cocoa Powder row).
# ---------------------------------------------------------------------------
$ws.Range("M17").Value = 0.0
$ws.Range("N17").Value = 0.0

# ---------------------------------------------------------------------------
# 4) Add three new ingredient rows (18, 19, 20).
# ---------------------------------------------------------------------------

# --- Row 18: BAR70 Dark Couverture (Chocolate) - reuses row 15/16 style ---
$ws.Range("A15").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Chocolate"

$ws.Range("B15").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "BAR70 Dark Couverture"

$ws.Range("C15:J15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 587.0
$ws.Range("D18").Value = 7.0
$ws.Range("E18").Value = 50.3
$ws.Range("F18").Value = 31.0
$ws.Range("G18").Value = 39.7
$ws.Range("H18").Value = 0.0
$ws.Range("I18").Value = 34.0
$ws.Range("J18").Value = 0.0

$ws.Range("M15").Copy()
$ws.Range("K18").PasteSpecial(-4122)
$ws.Range("K18").Value = 0.0

$ws.Range("H15").Copy()
$ws.Range("L18").PasteSpecial(-4122)
$ws.Range("L18").Value = 0.0

$ws.Range("M15").Copy()
$ws.Range("M18").PasteSpecial(-4122)
$ws.Range("M18").Value = 100.0

$ws.Range("M15").Copy()
$ws.Range("N18").PasteSpecial(-4122)
$ws.Range("N18").Value = 100.0

$ws.Range("C15").Copy()
$ws.Range("O18").PasteSpecial(-4122)
$ws.Range("O18").Value = 200.0

# --- Row 19: BAR38 Milk Couverture (Chocolate) - brand new bold style ---
$ws.Range("A2").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = "Chocolate"
$ws.Range("A19").Font.Bold = $true
$ws.Range("A19").Font.Size = 11
$ws.Range("A19").Font.Name = "Calibri"

$ws.Range("B15").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B19").Value = "BAR38 Milk Couverture"
$ws.Range("B19").Font.Bold = $true
$ws.Range("B19").Font.Size = 11
$ws.Range("B19").Font.Name = "Calibri"

$ws.Range("C15").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C19").Value = 560.0
$ws.Range("C19").Font.Size = 11
$ws.Range("C19").Font.Name = "Calibri"

$ws.Range("C15").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D19").Value = 4.2
$ws.Range("D19").Font.Size = 11
$ws.Range("D19").Font.Name = "Calibri"

$ws.Range("C15").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E19").Value = 56.0
$ws.Range("E19").Font.Size = 11
$ws.Range("E19").Font.Name = "Calibri"

$ws.Range("C15").Copy()
$ws.Range("F19").PasteSpecial(-4122)
$ws.Range("F19").Value = 47.0
$ws.Range("F19").Font.Size = 11
$ws.Range("F19").Font.Name = "Calibri"

$ws.Range("C15").Copy()
$ws.Range("G19").PasteSpecial(-4122)
$ws.Range("G19").Value = 35.0
$ws.Range("G19").Font.Size = 11
$ws.Range("G19").Font.Name = "Calibri"

$ws.Range("H15").Copy()
$ws.Range("H19").PasteSpecial(-4122)
$ws.Range("H19").Value = 0.0
$ws.Range("H19").Font.Size = 11
$ws.Range("H19").Font.Name = "Calibri"

$ws.Range("C15").Copy()
$ws.Range("I19").PasteSpecial(-4122)
$ws.Range("I19").Value = 0.0
$ws.Range("I19").Font.Size = 11
$ws.Range("I19").Font.Name = "Calibri"

$ws.Range("H15").Copy()
$ws.Range("J19").PasteSpecial(-4122)
$ws.Range("J19").Value = 0.0
$ws.Range("J19").Font.Size = 11
$ws.Range("J19").Font.Name = "Calibri"

$ws.Range("M15").Copy()
$ws.Range("K19").PasteSpecial(-4122)
$ws.Range("K19").Formula = "=1.07*(0.88*D19+6)"
$ws.Range("K19").Font.Size = 11
$ws.Range("K19").Font.Name = "Calibri"

$ws.Range("H15").Copy()
$ws.Range("L19").PasteSpecial(-4122)
$ws.Range("L19").Value = 0.0
$ws.Range("L19").Font.Size = 11
$ws.Range("L19").Font.Name = "Calibri"

$ws.Range("M15").Copy()
$ws.Range("M19").PasteSpecial(-4122)
$ws.Range("M19").Formula = "=((46*100)+(5*100))/51"
$ws.Range("M19").Font.Size = 11
$ws.Range("M19").Font.Name = "Calibri"

$ws.Range("N16").Copy()
$ws.Range("N19").PasteSpecial(-4122)
$ws.Range("N19").Formula = "=((46*100)+(5*16))/51"
$ws.Range("N19").Font.Size = 11
$ws.Range("N19").Font.Name = "Calibri"

$ws.Range("C15").Copy()
$ws.Range("O19").PasteSpecial(-4122)
$ws.Range("O19").Value = 145.0
$ws.Range("O19").Font.Size = 11
$ws.Range("O19").Font.Name = "Calibri"

# trailing formatted-but-empty cells P19:AC19
$ws.Range("H15").Copy()
$ws.Range("P19:AC19").PasteSpecial(-4122)
$ws.Range("P19:AC19").Font.Size = 11
$ws.Range("P19:AC19").Font.Name = "Calibri"
$ws.Range("P19:AC19").HorizontalAlignment = -4108
$ws.Range("P19:AC19").HorizontalAlignment = -4105

# --- Row 20: BAR Cocoa Powder (Cocoa Powder) - reuses row 17 style ---
$ws.Range("A17").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = "Cocoa Powder"

$ws.Range("B17").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B20").Value = "BAR Cocoa Powder"

$ws.Range("C17:O17").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 313.0
$ws.Range("D20").Value = 0.0
$ws.Range("E20").Value = 45.7
$ws.Range("F20").Value = 1.0
$ws.Range("G20").Value = 11.0
$ws.Range("H20").Value = 0.0
$ws.Range("I20").Value = 25.0
$ws.Range("J20").Value = 0.0
$ws.Range("K20").Value = 0.0
$ws.Range("L20").Value = 0.0
$ws.Range("M20").Value = 0.0
$ws.Range("N20").Value = 0.0
$ws.Range("O20").Value = 200.0

Write-Host "edit complete"
